$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Record results for Contest 26 (SRH vs RR, row 35) and Contest 27
#    (MI vs DC, row 36) -- these rows already had the match label, only the
#    per-player scores were missing.
# ---------------------------------------------------------------------------
$ws.Range("E35").Value = 100
$ws.Range("H35").Value = 60
$ws.Range("K35").Value = 40
$ws.Range("N35").Value = 0
$ws.Range("Q35").Value = 80
$ws.Range("T35").Value = 20

$ws.Range("E36").Value = 0
$ws.Range("H36").Value = 100
$ws.Range("K36").Value = 60
$ws.Range("N36").Value = 40
$ws.Range("Q36").Value = 80
$ws.Range("T36").Value = 20

# ---------------------------------------------------------------------------
# 2) Insert two new contest rows right before the old row 44 (the blank
#    placeholder row), shifting everything below down by two rows.
# ---------------------------------------------------------------------------
$ws.Range("A44:U45").Insert(-4121)

# The insert doesn't carry the row-43 styling into the new rows, so copy
# formats only (skip the 1-column spacers that have no style of their own,
# so we don't manufacture stray empty cells there).
$styledCols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T")
foreach ($col in $styledCols) {
    $ws.Range($col + "43").Copy()
    $ws.Range($col + "44:" + $col + "45").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Populate the two new contest rows: Contest 35 (SRH vs KKR) and
#    Contest 36 (MI vs KXI). Scores are still unknown, so only the
#    match-number/format/label and the standard ranking formulas go in.
# ---------------------------------------------------------------------------
$ws.Range("A44").Value = 35
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = "SRH vs KKR"

$ws.Range("A45").Value = 36
$ws.Range("B45").Value = 1
$ws.Range("C45").Value = "MI vs KXI"

$scoreToFormulaCol = @{ "D" = "E"; "G" = "H"; "J" = "K"; "M" = "N"; "P" = "Q"; "S" = "T" }
foreach ($r in 44, 45) {
    foreach ($fcol in $scoreToFormulaCol.Keys) {
        $scol = $scoreToFormulaCol[$fcol]
        $formula = '=IF(ISERROR(VLOOKUP(RANK(' + $scol + $r + ', ($T' + $r + ',$Q' + $r + ',$N' + $r + ',$K' + $r + ',$H' + $r + ',$E' + $r + '), 0),  score, 2, FALSE)),"",VLOOKUP(RANK(' + $scol + $r + ', ($T' + $r + ',$Q' + $r + ',$N' + $r + ',$K' + $r + ',$H' + $r + ',$E' + $r + '), 0),  score, 2, FALSE))'
        $ws.Range($fcol + $r).Formula = $formula
    }
}

# ---------------------------------------------------------------------------
# 4) The conditional formatting on the "Total" row still points at the old
#    row 48; move it to the new row 50 (this preserves dxfId/priority,
#    unlike deleting + re-adding the rules).
# ---------------------------------------------------------------------------
$totalCols = @("E", "H", "K", "N", "Q", "T")
foreach ($col in $totalCols) {
    $fcs = $ws.Range($col + "48").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($col + "50"))
    }
}

# ---------------------------------------------------------------------------
# 5) Selection follows the moved "Total" row.
# ---------------------------------------------------------------------------
$ws.Range("U50").Select()

$wb.Application.Calculate()
